$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("api_doc")

# F2: /api/health 500 error -- "false" -> "False"
$ws.Range("F2").Value = "500: { `"ok`": False, `"error`": `"internal_error`", “now”: …, “response_time_ms”: 120}"
$ws.Range("F2").Characters(1,3).Font.Bold = $true

# D3: /api/stats request description (unchanged text, refreshed)
$ws.Range("D3").Value = "Query params only (both mandatory)"

# E3: /api/stats success payload -- "status": "success" -> "ok": True; rest unchanged
$ws.Range("E3").Value = "{`"ok`": True, `"num_records`": 120, `"db_connected`": true, `"params`": {`"start_date`": …, `"end_date`": …}, `"response_time_ms`": 120, `"now`": …, `"data`": {`"day`":{`"ave`":{`"Fri`":396.25,`"Mon`":423.0,`"Sat`":360.25,`"Sun`":473.25,`"Thu`":355.25,`"Tue`":327.5,`"Wed`":357.25},`"std`":{`"Fri`":88.205,`"Mon`":44.728,`"Sat`":131.988,`"Sun`":120.477,`"Thu`":129.113,`"Tue`":147.789,`"Wed`":190.902}},`"week`":{`"ave`":2692.75,`"std`":261.545}}"

# F3: /api/stats 400 general_error payload -- drop "params", "status":"failure" -> "ok": True
$ws.Range("F3").Value = "400 (general_error): {`"ok`": True, `"error`": <python error str>, `"response_time_ms`": …, `"now`": …}"
$ws.Range("F3").Characters(1,3).Font.Bold = $true

# Update active selection to F3 (matches the authored cursor position)
$ws.Range("F3").Select() | Out-Null
